# Fruta / hortaliza, semanal
# Insert a new week of data (4 rows) at the top of the existing data block
# (before current row 26), shifting all subsequent rows down by 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows at rows 26-29, pushing existing data down.
$ws.Rows("26:29").Insert()

# Common values shared by every data row in this sheet.
$mercadoId = 1
$mercado   = "Agrícola del Norte S.A. de Arica"
$region    = "Arica y Parinacota"
$codreg    = 15
$tipo      = "Fruta"
$productoId = 100101
$producto   = "Berries"
$categoriaId = 100112025
$categoria   = "Frutilla"
$variedad    = "Sin especificar"
$unidad      = "`$/bandeja 3 kilos"
$origen      = "Región de Arica y Parinacota"
$kgUnidad    = 3

# New rows (fecha = 2022-11-30 -> serial 44895)
$fecha = 44895

$newRows = @(
    @{ Row = 26; Calidad = "Especial"; Volumen = 60; PMin = 7000; PMax = 8000; PProm = 7500; PKg = 2500 },
    @{ Row = 27; Calidad = "Primera";  Volumen = 70; PMin = 6000; PMax = 7000; PProm = 6500; PKg = 2167 },
    @{ Row = 28; Calidad = "Segunda";  Volumen = 72; PMin = 5000; PMax = 6000; PProm = 5486; PKg = 1829 },
    @{ Row = 29; Calidad = "Tercera";  Volumen = 74; PMin = 4000; PMax = 5000; PProm = 4500; PKg = 1500 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2  = $mercadoId
    $ws.Cells.Item($row, 2).Value2  = $mercado
    $ws.Cells.Item($row, 3).Value2  = $region
    $ws.Cells.Item($row, 4).Value2  = $fecha
    $ws.Cells.Item($row, 5).Value2  = $codreg
    $ws.Cells.Item($row, 6).Value2  = $tipo
    $ws.Cells.Item($row, 7).Value2  = $productoId
    $ws.Cells.Item($row, 8).Value2  = $producto
    $ws.Cells.Item($row, 9).Value2  = $categoriaId
    $ws.Cells.Item($row, 10).Value2 = $categoria
    $ws.Cells.Item($row, 11).Value2 = $variedad
    $ws.Cells.Item($row, 12).Value2 = $r.Calidad
    $ws.Cells.Item($row, 13).Value2 = $r.Volumen
    $ws.Cells.Item($row, 14).Value2 = $r.PMin
    $ws.Cells.Item($row, 15).Value2 = $r.PMax
    $ws.Cells.Item($row, 16).Value2 = $r.PProm
    $ws.Cells.Item($row, 17).Value2 = $unidad
    $ws.Cells.Item($row, 18).Value2 = $origen
    $ws.Cells.Item($row, 19).Value2 = $r.PKg
    $ws.Cells.Item($row, 20).Value2 = $kgUnidad
}
